$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some updated cells contain values that look like plain decimal numbers
# (e.g. "43.90", "0.619"). The source workbook stores every Price/Volume
# cell as text (inline string), so if we let Excel auto-detect these as
# numbers it would silently drop trailing zeros / switch to scientific
# notation (e.g. "43.90" -> 43.9, "0.0693" -> 6.93E-2). To avoid that we
# temporarily force a Text number format before writing the value, then
# restore the original (default/"Normal") style once the text is in place
# so no visible formatting change is introduced.
$textCells = @(
    "D5", "D6", "D8", "D9", "D10", "D13", "D15", "D16", "D18", "D20", "D22", "D26", "D27", "D29", "D32", "D33", "D36", "D37", "D38", "D39", "D41", "D43", "D44", "D47", "D48", "D51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "35.322.40"
$ws.Range("D3").Value = "1.842.74"
$ws.Range("E3").Value = "  +2.15%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Value = "233.97"
$ws.Range("E5").Value = "  +4.17%  "
$ws.Range("D6").Value = "0.619"
$ws.Range("E6").Value = "  +2.98%  "
$ws.Range("E7").Value = "  +0.34%  "
$ws.Range("D8").Value = "43.90"
$ws.Range("E8").Value = "  +11.53%  "
$ws.Range("D9").Value = "0.312"
$ws.Range("E9").Value = "  +7.66%  "
$ws.Range("D10").Value = "0.0693"
$ws.Range("E11").Value = "  +0.90%  "
$ws.Range("D12").Value = "2.108.78"
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").Value = "11.29"
$ws.Range("E13").Value = "  +3.57%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.830.13"
$ws.Range("E14").Value = "  +1.20%  "
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").Value = "0.671"
$ws.Range("E15").Value = "  +5.87%  "
$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D16").Value = "4.74"
$ws.Range("E16").Value = "  +8.71%  "
$ws.Range("D17").Value = "35.242.66"
$ws.Range("E17").Value = "  +1.75%  "
$ws.Range("D18").Value = "70.68"
$ws.Range("E18").Value = "  +4.34%  "
$ws.Range("D19").Value = "0.0₃0797"
$ws.Range("E19").Value = "  +3.98%  "
$ws.Range("D20").Value = "241.70"
$ws.Range("E20").Value = "  +0.82%  "
$ws.Range("E21").Value = "  +8.53%  "
$ws.Range("D22").Value = "4.63"
$ws.Range("E22").Value = "  +13.67%  "
$ws.Range("E23").Value = "  +0.52%  "
$ws.Range("E24").Value = "  +3.65%  "
$ws.Range("E25").Value = "  -0.07%  "
$ws.Range("D26").Value = "7.89"
$ws.Range("E26").Value = "  +3.16%  "
$ws.Range("D27").Value = "17.68"
$ws.Range("E27").Value = "  +1.69%  "
$ws.Range("E28").Value = "  +1.04%  "
$ws.Range("D29").Value = "1.58"
$ws.Range("E29").Value = "  +28.96%  "
$ws.Range("E30").Value = "  +0.82%  "
$ws.Range("D31").Value = "3.352.09"
$ws.Range("E31").Value = "  +37.96%  "
$ws.Range("D32").Value = "0.0558"
$ws.Range("E32").Value = "  +8.69%  "
$ws.Range("D33").Value = "3.95"
$ws.Range("E33").Value = "  +5.38%  "
$ws.Range("E34").Value = "  +6.28%  "
$ws.Range("E35").Value = "  +0.72%  "
$ws.Range("D36").Value = "94.86"
$ws.Range("E36").Value = "  +14.77%  "
$ws.Range("D37").Value = "0.687"
$ws.Range("E37").Value = "  +7.30%  "
$ws.Range("D38").Value = "1.11"
$ws.Range("E38").Value = "  +5.50%  "
$ws.Range("D39").Value = "0.0195"
$ws.Range("E39").Value = "  +3.92%  "
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").Value = "1.327.80"
$ws.Range("E40").Value = "  +1.87%  "
$ws.Range("B41").Value = "InjectiveProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D41").Value = "15.30"
$ws.Range("E41").Value = "  +4.56%  "
$ws.Range("E42").Value = "  +7.29%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "2.39"
$ws.Range("E43").Value = "  +1.36%  "
$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").Value = "1.28"
$ws.Range("E44").Value = "  +2.76%  "
$ws.Range("E45").Value = "  +0.36%  "
$ws.Range("E46").Value = "  -0.16%  "
$ws.Range("D47").Value = "6.24"
$ws.Range("E47").Value = "  +8.99%  "
$ws.Range("D48").Value = "0.0514"
$ws.Range("E48").Value = "  -0.67%  "
$ws.Range("D49").Value = "2.018.46"
$ws.Range("E49").Value = "  +2.59%  "
$ws.Range("E50").Value = "  +0.66%  "
$ws.Range("D51").Value = "101.59"
$ws.Range("E51").Value = "  +0.07%  "

foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}

